$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 1446.8889
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 1443.7646
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 4331.293799999999
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -6547.293799999999
# Row 132
$ws.Range("H132").Value = 56135.668
$ws.Range("I132").Value = 62971.375
$ws.Range("J132").Value = 1450
$ws.Range("K132").Value = 188914.125
$ws.Range("L132").Value = 4350
$ws.Range("M132").Value = -186384.125
$ws.Range("N132").Value = -9410
# Row 135
$ws.Range("H135").Value = 1430.9
$ws.Range("I135").Value = 1875.5
$ws.Range("J135").Value = 1134.5
$ws.Range("K135").Value = 16879.5
$ws.Range("L135").Value = 10210.5
$ws.Range("M135").Value = -14344.5
$ws.Range("N135").Value = -15280.5
# Row 137
$ws.Range("H137").Value = 1836.9474
$ws.Range("I137").Value = 2295.1
$ws.Range("K137").Value = 6885.299999999999
$ws.Range("M137").Value = -4335.299999999999
# Row 138
$ws.Range("H138").Value = 15387057
$ws.Range("I138").Value = 1177.1562
$ws.Range("J138").Value = 30306698
$ws.Range("K138").Value = 3531.4686
$ws.Range("L138").Value = 90920094
$ws.Range("M138").Value = 1608.5314
$ws.Range("N138").Value = -90930374

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 18
$ws.Range("H18").Value = 3506.5
$ws.Range("J18").Value = 3506.5
$ws.Range("L18").Value = 3506.5
$ws.Range("N18").Value = -4150.5
# Row 122
$ws.Range("H122").Value = 1700.069
$ws.Range("I122").Value = 2026.5454
$ws.Range("J122").Value = 1500.5555
$ws.Range("K122").Value = 6079.6362
$ws.Range("L122").Value = 4501.666499999999
$ws.Range("M122").Value = -3629.6362
$ws.Range("N122").Value = -9401.666499999999
# Row 132
$ws.Range("H132").Value = 2643.3635
$ws.Range("I132").Value = 1022.1053
$ws.Range("K132").Value = 3066.3159
$ws.Range("M132").Value = -536.3159000000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 20005418
$ws.Range("I31").Value = 40000870
$ws.Range("J31").Value = 9966.120000000001
$ws.Range("K31").Value = 40000870
$ws.Range("L31").Value = 9966.120000000001
$ws.Range("M31").Value = -40000575
$ws.Range("N31").Value = -10556.12
# Row 34
$ws.Range("H34").Value = 20005418
$ws.Range("I34").Value = 40000870
$ws.Range("J34").Value = 9966.120000000001
$ws.Range("K34").Value = 40000870
$ws.Range("L34").Value = 9966.120000000001
$ws.Range("M34").Value = -40000668
$ws.Range("N34").Value = -10370.12
# Row 99
$ws.Range("H99").Value = 15444.444
$ws.Range("I99").Value = 5314.2856
$ws.Range("J99").Value = 50900
$ws.Range("K99").Value = 5314.2856
$ws.Range("L99").Value = 50900
$ws.Range("M99").Value = -3816.2856
$ws.Range("N99").Value = -53896
# Row 126
$ws.Range("H126").Value = 15444.444
$ws.Range("I126").Value = 5314.2856
$ws.Range("J126").Value = 50900
$ws.Range("K126").Value = 15942.8568
$ws.Range("L126").Value = 152700
$ws.Range("M126").Value = -13472.8568
$ws.Range("N126").Value = -157640
# Row 134
$ws.Range("H134").Value = 2891.46
$ws.Range("I134").Value = 3721.3872
$ws.Range("K134").Value = 11164.1616
$ws.Range("M134").Value = -8629.161599999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1239383.8
$ws.Range("J68").Value = 1135.24
$ws.Range("L68").Value = 3405.72
$ws.Range("N68").Value = -5027.72
# Row 71
$ws.Range("H71").Value = 1239383.8
$ws.Range("J71").Value = 1135.24
$ws.Range("L71").Value = 10217.16
$ws.Range("N71").Value = -18329.16
# Row 112
$ws.Range("H112").Value = 1785
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1785
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -7571
# Row 113
$ws.Range("H113").Value = 1717.52
$ws.Range("I113").Value = 1759.4
$ws.Range("J113").Value = 1654.7
$ws.Range("K113").Value = 5278.200000000001
$ws.Range("L113").Value = 4964.1
$ws.Range("M113").Value = -3108.200000000001
$ws.Range("N113").Value = -9304.1
# Row 136
$ws.Range("H136").Value = 1515.909
$ws.Range("I136").Value = 1515.909
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4547.727000000001
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2799.8
$ws.Range("I102").Value = 2749.75
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 2749.75
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -1127.75
$ws.Range("N102").Value = -6244
# Row 132
$ws.Range("H132").Value = 3634.7144
$ws.Range("I132").Value = 1654
$ws.Range("J132").Value = 4427
$ws.Range("K132").Value = 4962
$ws.Range("L132").Value = 13281
$ws.Range("M132").Value = -2432
$ws.Range("N132").Value = -18341

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1898.6072
$ws.Range("I16").Value = 1020.8333
$ws.Range("J16").Value = 3478.6
$ws.Range("K16").Value = 1020.8333
$ws.Range("L16").Value = 3478.6
$ws.Range("M16").Value = -850.8333
$ws.Range("N16").Value = -3818.6
# Row 132
$ws.Range("H132").Value = 6737.727
$ws.Range("I132").Value = 8926.799999999999
$ws.Range("J132").Value = 4913.5
$ws.Range("K132").Value = 26780.4
$ws.Range("L132").Value = 14740.5
$ws.Range("M132").Value = -24250.4
$ws.Range("N132").Value = -19800.5
# Row 136
$ws.Range("H136").Value = 10103074
$ws.Range("I136").Value = 1049.0625
$ws.Range("J136").Value = 19610862
$ws.Range("K136").Value = 3147.1875
$ws.Range("L136").Value = 58832586
$ws.Range("M136").Value = -597.1875
$ws.Range("N136").Value = -58837686

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 64964.4
$ws.Range("I20").Value = 57400
$ws.Range("J20").Value = 70007.336
$ws.Range("K20").Value = 57400
$ws.Range("L20").Value = 70007.336
$ws.Range("M20").Value = -57160
$ws.Range("N20").Value = -70487.336
# Row 132
$ws.Range("H132").Value = 1836.4286
$ws.Range("I132").Value = 939.75
$ws.Range("J132").Value = 3032
$ws.Range("K132").Value = 2819.25
$ws.Range("L132").Value = 9096
$ws.Range("M132").Value = -289.25
$ws.Range("N132").Value = -14156
# Row 136
$ws.Range("H136").Value = 4539.3613
$ws.Range("I136").Value = 670.41174
$ws.Range("K136").Value = 2011.23522
$ws.Range("M136").Value = 538.76478
